# ProjectApplication.xlsx update
# - adds a "Flat Type" column between "Application Status" and "Date"
# - updates the sample application record values
# - changes the date-formatted style to a 2-decimal numeric style
# - updates sheet view (zoom/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column for "Flat Type" before the existing "Date" column (E),
# shifting Date to column F.
$ws.Columns("E:E").Insert()

# ---- Header row ----
$ws.Cells.Item(1,1).Value = "Application ID"
$ws.Cells.Item(1,2).Value = "Project ID"
$ws.Cells.Item(1,3).Value = "Applicant NRIC"
$ws.Cells.Item(1,4).Value = "Application Status"
$ws.Cells.Item(1,5).Value = "Flat Type"
$ws.Cells.Item(1,6).Value = "Date"

# The style previously used for the date-formatted cell (now F2, shifted from
# the old E2) is repurposed in-place as a 2-decimal numeric style, and then
# reused on the "Application ID"/"Project ID" headers.
$ws.Range("F2").NumberFormat = "0.00"
$ws.Range("A1").NumberFormat = "0.00"
$ws.Range("B1").NumberFormat = "0.00"

# ---- Data row ----
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = "S1234567A"
$ws.Cells.Item(2,4).Value = "Withdrawal Pending"
$ws.Cells.Item(2,5).Value = "2-ROOM"
$ws.Cells.Item(2,6).Value = 45767.18201704861

# ---- View settings ----
$excel.ActiveWindow.Zoom = 85
$ws.Range("A3:XFD1048576").Select() | Out-Null

Write-Host "ProjectApplication.xlsx updated"
